$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values that look numeric (e.g. "64.155.29", "0.0000153")
# but must remain plain text, matching the original inlineStr cells. Force the
# whole Price column to Text format before writing the values so Excel does not
# reinterpret them as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.155.29"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "3.131.84"
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "568.79"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").Value = "161.18"
$ws.Range("E6").Value = "  -4.46%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "0.562"
$ws.Range("E8").Value = "  -7.63%  "
$ws.Range("E9").Value = "  -3.59%  "
$ws.Range("D10").Value = "6.56"
$ws.Range("E10").Value = "  -2.55%  "
$ws.Range("D11").Value = "0.379"
$ws.Range("E11").Value = "  -1.60%  "
$ws.Range("D12").Value = "3.682.27"
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("D14").Value = "64.285.24"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").Value = "24.83"
$ws.Range("E15").Value = "  -1.88%  "
$ws.Range("D16").Value = "3.142.25"
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("D17").Value = "0.0000153"
$ws.Range("E17").Value = "  -3.11%  "
$ws.Range("D18").Value = "401.09"
$ws.Range("E18").Value = "  -4.46%  "
$ws.Range("D19").Value = "12.48"
$ws.Range("E19").Value = "  -2.56%  "
$ws.Range("D20").Value = "5.21"
$ws.Range("E20").Value = "  -2.65%  "
$ws.Range("D21").Value = "7.09"
$ws.Range("E21").Value = "  +0.79%  "
$ws.Range("E22").Value = "  +3.93%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").Value = "67.82"
$ws.Range("E24").Value = "  -3.27%  "
$ws.Range("D25").Value = "0.479"
$ws.Range("E25").Value = "  -1.47%  "
$ws.Range("E26").Value = "  -4.21%  "
$ws.Range("D27").Value = "0.0000100"
$ws.Range("E27").Value = "  -4.83%  "
$ws.Range("D28").Value = "8.76"
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.80"
$ws.Range("E31").Value = "  -1.38%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "21.12"
$ws.Range("E32").Value = "  -2.94%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "158.72"
$ws.Range("E33").Value = "  +1.00%  "
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").Value = "6.22"
$ws.Range("E34").Value = "  -1.36%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "4.79"
$ws.Range("E35").Value = "  -4.50%  "
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").Value = "1.10"
$ws.Range("E36").Value = "  -2.84%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "1.33"
$ws.Range("E37").Value = "  -2.37%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "2.657.05"
$ws.Range("E38").Value = "  -1.63%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "1.66"
$ws.Range("E39").Value = "  -2.07%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "23.41"
$ws.Range("E40").Value = "  -3.72%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "4.06"
$ws.Range("E41").Value = "  -2.22%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "38.29"
$ws.Range("E42").Value = "  -2.17%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "0.686"
$ws.Range("E43").Value = "  -3.20%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").Value = "0.0608"
$ws.Range("E44").Value = "  -1.56%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "5.45"
$ws.Range("E45").Value = "  -2.86%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "0.0254"
$ws.Range("E46").Value = "  -2.77%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").Value = "286.80"
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "20.90"
$ws.Range("E48").Value = "  -3.34%  "
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").Value = "0.997"
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "0.0970"
$ws.Range("E50").Value = "  -1.80%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "10.48"
$ws.Range("E51").Value = "  +0.57%  "
